$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the time values in B7:F7 while keeping their formatting
$ws.Range("B7:F7").ClearContents()

# Update the active selection to F7
$ws.Range("F7").Select()
